$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 11 de Agosto de 2020 a las 15:31"
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 5253765
$ws.Range("C4").Value = 2319
$ws.Range("D4").Value = 2717257
$ws.Range("E4").Value = 2370235
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 81
$ws.Range("H4").Value = 166273

$ws.Range("A6").Value = "India"
$ws.Range("B6").Value = 2287511
$ws.Range("C6").Value = 20358
$ws.Range("D6").Value = 1598529
$ws.Range("E6").Value = 643432
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 197
$ws.Range("H6").Value = 45550

$ws.Range("A16").Value = "Arabia Saudita"
$ws.Range("B16").Value = 291468
$ws.Range("C16").Value = 1521
$ws.Range("D16").Value = 255118
$ws.Range("E16").Value = 33117
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 34
$ws.Range("H16").Value = 3233

$ws.Range("A24").Value = "Irak"
$ws.Range("B24").Value = 156995
$ws.Range("C24").Value = 3396
$ws.Range("D24").Value = 112102
$ws.Range("E24").Value = 39362
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 67
$ws.Range("H24").Value = 5531

$ws.Range("A35").Value = "Suecia"
$ws.Range("B35").Value = 83126
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 5
$ws.Range("H35").Value = 5770

$ws.Range("A36").Value = "Ucrania"
$ws.Range("B36").Value = 83115
$ws.Range("C36").Value = 1158
$ws.Range("D36").Value = 44934
$ws.Range("E36").Value = 36230
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 29
$ws.Range("H36").Value = 1951

$ws.Range("A45").Value = "Paises Bajos"
$ws.Range("B45").Value = 59973
$ws.Range("C45").Value = 779
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = 6159

$ws.Range("A53").Value = "Barein"
$ws.Range("B53").Value = 44397
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 41209
$ws.Range("E53").Value = 3023
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 2
$ws.Range("H53").Value = 165

$ws.Range("A61").Value = "Azerbaiyan"
$ws.Range("B61").Value = 33731
$ws.Range("C61").Value = 84
$ws.Range("D61").Value = 30856
$ws.Range("E61").Value = 2380
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 3
$ws.Range("H61").Value = 495

$ws.Range("A62").Value = "Uzbekistan"
$ws.Range("B62").Value = 31747
$ws.Range("C62").Value = 443
$ws.Range("D62").Value = 23704
$ws.Range("E62").Value = 7839
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 4
$ws.Range("H62").Value = 204

$ws.Range("A63").Value = "Serbia"
$ws.Range("B63").Value = 28497
$ws.Range("C63").Value = 235
$ws.Range("D63").Value = 18965
$ws.Range("E63").Value = 8880
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 6
$ws.Range("H63").Value = 652

$ws.Range("A65").Value = "Kenia"
$ws.Range("B65").Value = 27425
$ws.Range("C65").Value = 497
$ws.Range("D65").Value = 13867
$ws.Range("E65").Value = 13120
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 15
$ws.Range("H65").Value = 438

$ws.Range("A78").Value = "Estado de Palestina"
$ws.Range("B78").Value = 14875
$ws.Range("C78").Value = 365
$ws.Range("D78").Value = 8181
$ws.Range("E78").Value = 6591
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 3
$ws.Range("H78").Value = 103

$ws.Range("A79").Value = "Bosnia y Herzegovina"
$ws.Range("B79").Value = 14708
$ws.Range("C79").Value = 210
$ws.Range("D79").Value = 8411
$ws.Range("E79").Value = 5850
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 22
$ws.Range("H79").Value = 447

$ws.Range("A80").Value = "Corea del Sur"
$ws.Range("B80").Value = 14660
$ws.Range("C80").Value = 34
$ws.Range("D80").Value = 13729
$ws.Range("E80").Value = 626
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 305

$ws.Range("A83").Value = "Republica de Macedonia"
$ws.Range("B83").Value = 12083
$ws.Range("C83").Value = 141
$ws.Range("D83").Value = 8248
$ws.Range("E83").Value = 3306
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = 529

$ws.Range("A84").Value = "Sudan"
$ws.Range("B84").Value = 12033
$ws.Range("C84").Value = 77
$ws.Range("D84").Value = 6282
$ws.Range("E84").Value = 4965
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 5
$ws.Range("H84").Value = 786

$ws.Range("A103").Value = "Croacia"
$ws.Range("B103").Value = 5740
$ws.Range("C103").Value = 91
$ws.Range("D103").Value = 4962
$ws.Range("E103").Value = 618
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 2
$ws.Range("H103").Value = 160

$ws.Range("A134").Value = "Islandia"
$ws.Range("B134").Value = 1968
$ws.Range("C134").Value = 6
$ws.Range("D134").Value = 1844
$ws.Range("E134").Value = 114
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 10

$ws.Range("A159").Value = "Vietnam"
$ws.Range("B159").Value = 863
$ws.Range("C159").Value = 16
$ws.Range("D159").Value = 399
$ws.Range("E159").Value = 448
$ws.Range("F159").Value = 0
$ws.Range("G159").Value = 1
$ws.Range("H159").Value = 16

$ws.Range("A213").Value = "Montserrat"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 12
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

